$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The published table grew a new "2020" column (Q), mirroring the style of
# the existing "2019" column (P) for the header (row 3) and the data value
# (row 4).

# 1) Change the custom "0.0" number format used by the data row (row 4,
#    columns D:P) to the built-in "0.00" format before extending the row,
#    so the new Q4 cell picks up the same (already-updated) formatting
#    when it is copied from P4 below.
$ws.Range("D4:P4").NumberFormat = "0.00"

# 2) Add the new year column: copy formatting from the last existing
#    column (P) into the new column (Q) for both the header row and the
#    data row, then overwrite with the new values.
$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2020

$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 0.067156049127444606

# 3) Reset the selection back to the top-left cell (the sheet was left
#    with a stray selection on B12 from a prior editing session).
$null = $ws.Range("A1").Select()
